$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, shifting the existing row 20 (and everything
# below it) down by one. Excel's Insert() copies formatting from the row
# above into the freshly inserted row, so we clear that spillover back out
# to leave row 20 with only its (pre-existing) A20 styling.
$ws.Rows.Item(20).Insert()
$ws.Range("B20:D20").Clear()

# Match the new selection recorded in the workbook.
$ws.Range("B20").Select()
